$wb = $excel.ActiveWorkbook
$todo = $wb.Worksheets.Item("To Do")
$notif = $wb.Worksheets.Item("Notifications")

# --- sheet "Notifications": update B2 text (Observer registered for User too) ---
$notif.Range("B2").Value = "Donator - Creator - Admin - User"

# --- sheet "To Do": add two new task rows (31 & 32) ---
$todo.Range("A31").Value = "Notification Against Comment"
$todo.Range("B31").Value = "++"
$todo.Range("A32").Value = "Email Against Transaction"
$todo.Range("B32").Value = "++"

# match formatting of the other "in progress" rows (B column quote-prefixed style)
$todo.Range("B29").Copy()
$todo.Range("B31:B32").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- sheet "Notifications": new row 11 (comment notification) ---
$notif.Range("A11").Value = "on Comments"

# restore view/selection to match the edited areas
$todo.Activate()
[void]$todo.Range("A31").Select()

$notif.Activate()
[void]$notif.Range("A11").Select()
